# Revert "Merge branch 'main' of https://github.com/kazarach/mitraaccounting"
# This reverts the merge commit, restoring Sheet1 of "FE JUAN.xlsx" to its
# pre-merge state:
#   - the highlighted "Filter Modal : category, supplier." note is renamed
#     back to "Modal : category, supplier." and re-inserted right after the
#     "Pesanan Penjualan" row instead of living near the bottom of the sheet
#   - the extra helper cells in column A (rows 3-6) and the accent fill that
#     had been applied across A3:C6 are removed again
#   - the sheet's active selection goes back to P17

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content: restore the original wording of the "category, supplier" note
# (previously "Filter Modal : category, supplier. ", now just
# "Modal : category, supplier. ") that lives in C5, next to "Pesanan Penjualan".
$ws.Range("C5").Value = "Modal : category, supplier. "

# --- Formatting: drop the accent-color fill that used to highlight A3:C6.
# Column A (B3:A6) only ever held the fill with no text, so those cells are
# removed outright; B3:C6 keep their text but lose the applied fill/style.
$ws.Range("A3:A6").Clear()
$ws.Range("B3:C6").ClearFormats()

# --- View state: put the active selection back on P17 (it had drifted to C10).
$ws.Range("P17").Select()
